$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.589.01"
$ws.Range("E2").Value = "  +5.01%  "
$ws.Range("D3").Value = "2.726.91"
$ws.Range("E3").Value = "  +3.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.65"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.47"
$ws.Range("E6").Value = "  +5.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("D9").Value = "2.754.66"
$ws.Range("E9").Value = "  +4.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.74"
$ws.Range("E10").Value = "  +3.51%  "
$ws.Range("E11").Value = "  +6.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.391"
$ws.Range("E12").Value = "  +4.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.161"
$ws.Range("E13").Value = "  +4.06%  "
$ws.Range("D14").Value = "3.229.40"
$ws.Range("E14").Value = "  +4.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.35"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "63.673.86"
$ws.Range("E16").Value = "  +5.16%  "
$ws.Range("E17").Value = "  +7.40%  "
$ws.Range("D18").Value = "2.750.86"
$ws.Range("E18").Value = "  +4.67%  "
$ws.Range("E19").Value = "  +4.09%  "
$ws.Range("E20").Value = "  +3.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "360.82"
$ws.Range("E21").Value = "  +3.46%  "
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.537"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.09"
$ws.Range("E25").Value = "  +3.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.168"
$ws.Range("E26").Value = "  +4.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.58"
$ws.Range("E27").Value = "  +5.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  +11.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.02"
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("E31").Value = "  +7.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.95"
$ws.Range("E32").Value = "  +2.25%  "
$ws.Range("E33").Value = "  +16.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.54"
$ws.Range("E35").Value = "  +4.90%  "
$ws.Range("E36").Value = "  +7.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.45"
$ws.Range("E37").Value = "  +9.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.82"
$ws.Range("E38").Value = "  +9.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  +15.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "344.00"
$ws.Range("E40").Value = "  +4.07%  "
$ws.Range("E41").Value = "  +5.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.13"
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.51"
$ws.Range("E43").Value = "  +7.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.70"
$ws.Range("E44").Value = "  +7.82%  "
$ws.Range("E45").Value = "  +7.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "139.61"
$ws.Range("E46").Value = "  +4.31%  "
$ws.Range("E47").Value = "  +5.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0587"
$ws.Range("E48").Value = "  +5.47%  "
$ws.Range("E49").Value = "  +5.07%  "
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  +0.02%  "
